# Insert a new "Currency" column before the existing "Invoice Value" column.
# This shifts the old E:H columns (Invoice Value, Unnamed: 5, Due, Ageing (Days))
# one place to the right, becoming F:I, and the new column E is populated with
# the currency symbol implied by each row's invoice number series:
#   "...ES..." (export services, foreign clients)  -> "$"
#   "...LS..." (local services, domestic clients)  -> "₹"
# Row 2 (the opening-balance row) has no invoice number but belongs to the
# same $-denominated client block as the rows following it, so it defaults to "$".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E; existing E/F/G/H shift right to F/G/H/I.
$ws.Columns("E:E").Insert()

# Header row.
$ws.Range("E1").Value2 = "Currency"
$ws.Range("G1").Value2 = "Unnamed: 6"

# Fill the new Currency column for every data row (2-103).
$lastRow = 103
for ($r = 2; $r -le $lastRow; $r++) {
    $invoiceNo = $ws.Cells.Item($r, 4).Value2
    if ($invoiceNo -eq $null -or $invoiceNo -eq "") {
        $currency = "$"
    } elseif ($invoiceNo -like "*LS*") {
        $currency = "₹"
    } else {
        $currency = "$"
    }
    $ws.Cells.Item($r, 5).Value2 = $currency
}
